# Update sexting escalation script lines on the JockurworldJourney sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JockurworldJourney")

$ws.Range("B4").Value  = "watch me cum for you right now bro 💦"
$ws.Range("B6").Value  = "I'm cumming... FUCK watch every drop"
$ws.Range("B7").Value  = "I'm right there bro... don't stop watching, I'm about to blow 💦"
$ws.Range("B8").Value  = "holy fuck 🥵"
$ws.Range("B9").Value  = "you're about to see what happens when I completely let go 💦"
$ws.Range("B11").Value = "you need to watch what happens next bro... I'm about to explode"
$ws.Range("B12").Value = "I keep imagining you here and it's making me lose my mind 💦"
$ws.Range("B13").Value = "I'm going so hard right now and I can feel every stroke building bro... I'm dripping everywhere"
$ws.Range("B14").Value = "FUCK 🥵"
$ws.Range("B15").Value = "this is what you do to me bro... watch 💦"
$ws.Range("B17").Value = "tell me what you'd do if you were here right now... don't hold back"
$ws.Range("B18").Value = "I'm stroking myself thinking about you right now and I can barely handle it 💦"
$ws.Range("B19").Value = "I can't stop now... I'm throbbing so hard and pre-cum is already dripping"
$ws.Range("B20").Value = "damn bro... okay that was intense 🥵"
$ws.Range("B21").Value = "look what you did bro... you're not ready for this 💦"
$ws.Range("B23").Value = "I'm already gripping my cock and stroking it because of you... hope you can handle what you started 💦"
$ws.Range("B24").Value = "I can feel myself getting bigger just from talking to you... my body doesn't lie bro"
$ws.Range("B25").Value = "you liked that huh? because I'm getting hard just knowing you're looking 🥵"
